$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.030979633331299
$ws.Range("B1").Value = 0.9703243970870972
$ws.Range("C1").Value = 0.754185676574707
$ws.Range("D1").Value = 0.733465850353241
$ws.Range("E1").Value = 0.7948825359344482
